# "Agrego dato del IPC INDEC Oct24 y res fiscal del mismo mes"
#
# - IPC-Seriemensual: append the October-2024 monthly IPC row.
# - Aperturas / IPC-DIC-Div / IPC-Interanual: refresh the rolling 12-month
#   windows that shift forward once the new month lands.
# - Leave IPC-Seriemensual as the active/selected sheet when done (it was
#   EMAE before the edit).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) IPC-Seriemensual: add the new monthly row (Oct-2024)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("IPC-Seriemensual")

# Copy A94's format down to A95 first so the new date cell keeps the same
# date-number-format style as the rest of column A, then write the value.
$ws1.Range("A94").Copy()
$ws1.Range("A95").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("A95").Value = 45566
$ws1.Range("B95").Value = 2.7
$ws1.Range("C95").Value = 1.4
$ws1.Range("D95").Value = 2.9
$ws1.Range("E95").Value = 2.7

# ---------------------------------------------------------------------------
# 2) Aperturas: refreshed rolling figures (rows 2-10, cols B:H)
# ---------------------------------------------------------------------------
$ws12 = $wb.Worksheets.Item("Aperturas")
$aperturas = @(
  @(2.1, 2.2000000000000002, 1.2, 1.5, 0.7, 2.4, 1.9981262523715237),
  @(0.9, 1.2, 2.1, 1, 0.5, 1.5, 1.0692085999603806),
  @(2.5, 1.8, 1.7, 2.1, 2.1, 3.7, 2.2415449669569432),
  @(1.7, 1.9, 2.1, 2.2999999999999998, 3.3, 2, 1.9443776060163787),
  @(3.7, 6.4, 3.1, 6, 4.8, 4.2, 4.8629044308175828),
  @(-5, -3, -3.9, -4.3, -9.9, -3.3, -4.4659146924120758),
  @(4.2, 2.8, 2.4, 1.7, 2, 3, 3.2891960962084354),
  @(2.1, 2.1, 1.6, 1.9, 2.2000000000000002, -0.4, 1.9405868017069228),
  @(4.4000000000000004, 3.2, 2.4, 1.7, 2.4, 3, 3.4901795463939012)
)
for ($i = 0; $i -lt $aperturas.Length; $i++) {
  $r = $i + 2
  $row = $aperturas[$i]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws12.Cells.Item($r, $c + 2).Value = $row[$c]
  }
}

# ---------------------------------------------------------------------------
# 3) IPC-DIC-Div: refreshed rolling 12-month table (rows 2-14, cols B:H)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("IPC-DIC-Div")
$dicdiv = @(
  @(2.7, 2.8, 2.6, 2.6, 2.6, 2.2000000000000002, 3.2),
  @(1.2, 1.2, 1.5, 1.1000000000000001, 0.8, -0.4, 1.7),
  @(3, 3.6, 2.6, 2.1, 2.5, 2.4, 2.4),
  @(4.4000000000000004, 4.5, 4.3, 5, 4, 4.4000000000000004, 5.2),
  @(5.4, 4.7, 5.7, 7.2, 5, 6.7, 7.5),
  @(2.6, 2.6, 2.7, 2.2999999999999998, 2.6, 2.2000000000000002, 2.5),
  @(3.6, 3.9, 3.2, 3.1, 3, 3, 5.7),
  @(1.2, 1.5, 0.8, 0.6, 1.9, 1.7, 0.9),
  @(2.1, 2.2000000000000002, 2.1, 1.5, 1.3, 2, 2.6),
  @(2.9, 3.4, 2, 2.6, 3, 3.6, 3),
  @(3.5, 2.7, 3.9, 1.4, 7.4, 4.0999999999999996, 4),
  @(4.3, 4, 4, 5.0999999999999996, 5.6, 4.2, 6.2),
  @(2.8, 3.3, 2.5, 2.8, 2.1, 1.8, 3.1)
)
for ($i = 0; $i -lt $dicdiv.Length; $i++) {
  $r = $i + 2
  $row = $dicdiv[$i]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws2.Cells.Item($r, $c + 2).Value = $row[$c]
  }
}

# ---------------------------------------------------------------------------
# 4) IPC-Interanual: refreshed rolling 12-month table (rows 2-14, cols B:H)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("IPC-Interanual")
$interanual = @(
  @(193, 197.3, 190.4, 181.8, 184.5, 193.8, 195),
  @(183.2, 187.3, 183.9, 161.5, 168.3, 175.3, 193.1),
  @(184.1, 191.5, 182.2, 162.30000000000001, 168.2, 172.5, 189.4),
  @(130, 127.6, 132.19999999999999, 134.19999999999999, 122.5, 137.19999999999999, 139.80000000000001),
  @(285.8, 269.5, 284.5, 304.89999999999998, 344.8, 334.2, 290.10000000000002),
  @(166.2, 160, 174, 155.4, 166.2, 173.4, 173.6),
  @(220.2, 227, 216.3, 213.3, 210.3, 212.4, 221.4),
  @(226.8, 248.7, 207.8, 219.4, 230.7, 223.2, 194.4),
  @(257.60000000000002, 263.2, 254.8, 262.10000000000002, 235.4, 258.8, 254.4),
  @(170.7, 174.2, 170.2, 175.1, 162.80000000000001, 162, 156.1),
  @(188.5, 191.5, 173, 193.9, 214.2, 190.3, 218.6),
  @(184.4, 187.2, 181.4, 170.4, 179.9, 185.8, 201.8),
  @(247.5, 237.9, 254.4, 237.1, 246.8, 262.60000000000002, 289.3)
)
for ($i = 0; $i -lt $interanual.Length; $i++) {
  $r = $i + 2
  $row = $interanual[$i]
  for ($c = 0; $c -lt $row.Length; $c++) {
    $ws3.Cells.Item($r, $c + 2).Value = $row[$c]
  }
}

# ---------------------------------------------------------------------------
# 5) View / selection bookkeeping to match the author's final saved state
# ---------------------------------------------------------------------------
# IPC-Interanual: selection moves off the old H2:H14 block back to default A1.
$ws3.Range("A1").Select()

# IPC-DIC-Div: selection widens from H2:H14 to B2:H14.
$ws2.Range("B2:H14").Select()

# IPC-Seriemensual becomes the active sheet (was EMAE), selection on the new row.
$ws1.Range("C95:E95").Select()
$ws1.Activate()
